# Generate Report for Handoff
#
# Re-stamp the "Latest Handoff Datetime" (column D) for the files that were
# just (re-)handed off, on both locale status sheets. All of the affected
# rows receive the same timestamp per locale, since the handoff report was
# generated in a single batch run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnTimestamp = "2016-03-09 20:25:04"
foreach ($row in 7,10,11,12,13,14,15,16) {
    $zhcn.Range("D$row").Value = $zhcnTimestamp
}

$dede = $wb.Worksheets.Item("de-de")
$dedeTimestamp = "2016-03-09 20:25:10"
foreach ($row in 7,10,11,12,13,14,15,16) {
    $dede.Range("D$row").Value = $dedeTimestamp
}
